$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("E:E").Insert()
$ws.Columns("E").ColumnWidth = 29.329
